$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PO Forecast"

# --- Match the sheet-level layout of the existing sheets (outline +
#     page margins), mirroring <sheetPr><outlinePr .../></sheetPr> and
#     <pageMargins .../> on "Weekly Quantity" / "Monthly Trend". ---
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# --- Reuse the existing header style (bold + border + centered) from the
#     "Weekly Quantity" sheet, and the existing date-serial number format
#     style from its A column, so no new styles are introduced. ---
$wsWeekly.Range("A1:B1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$wsWeekly.Range("A2").Copy()
$ws.Range("A2:A30").PasteSpecial(-4122)

# --- Header row ---
$ws.Cells.Item(1,1).Value = "ds"
$ws.Cells.Item(1,2).Value = "PO_Forecast"
$ws.Cells.Item(1,3).Value = "yhat_lower"
$ws.Cells.Item(1,4).Value = "yhat_upper"

# --- Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$ws.Cells.Item(2,1).Value = 44934.99999999999
$ws.Cells.Item(2,2).Value = 40
$ws.Cells.Item(2,3).Value = -7.232265561684255
$ws.Cells.Item(2,4).Value = 85.06719518781991
$ws.Cells.Item(3,1).Value = 44962.99999999999
$ws.Cells.Item(3,2).Value = 42
$ws.Cells.Item(3,3).Value = -3.437698771952238
$ws.Cells.Item(3,4).Value = 86.88996360811369
$ws.Cells.Item(4,1).Value = 44969.99999999999
$ws.Cells.Item(4,2).Value = 42
$ws.Cells.Item(4,3).Value = -3.136084444214162
$ws.Cells.Item(4,4).Value = 86.48575293960864
$ws.Cells.Item(5,1).Value = 44976.99999999999
$ws.Cells.Item(5,2).Value = 43
$ws.Cells.Item(5,3).Value = -3.222396190101977
$ws.Cells.Item(5,4).Value = 88.8654117022361
$ws.Cells.Item(6,1).Value = 44983.99999999999
$ws.Cells.Item(6,2).Value = 43
$ws.Cells.Item(6,3).Value = -4.87988594232902
$ws.Cells.Item(6,4).Value = 87.42931962977497
$ws.Cells.Item(7,1).Value = 44990.99999999999
$ws.Cells.Item(7,2).Value = 44
$ws.Cells.Item(7,3).Value = -0.4624403581518806
$ws.Cells.Item(7,4).Value = 93.12760801056781
$ws.Cells.Item(8,1).Value = 44997.99999999999
$ws.Cells.Item(8,2).Value = 44
$ws.Cells.Item(8,3).Value = -5.306191891070566
$ws.Cells.Item(8,4).Value = 89.19145783982746
$ws.Cells.Item(9,1).Value = 45004.99999999999
$ws.Cells.Item(9,2).Value = 45
$ws.Cells.Item(9,3).Value = -2.107314606384973
$ws.Cells.Item(9,4).Value = 87.72508425555175
$ws.Cells.Item(10,1).Value = 45011.99999999999
$ws.Cells.Item(10,2).Value = 45
$ws.Cells.Item(10,3).Value = -3.004460618458441
$ws.Cells.Item(10,4).Value = 88.7679640292678
$ws.Cells.Item(11,1).Value = 45018.99999999999
$ws.Cells.Item(11,2).Value = 45
$ws.Cells.Item(11,3).Value = 0.1778452913631731
$ws.Cells.Item(11,4).Value = 90.35957315563715
$ws.Cells.Item(12,1).Value = 45025.99999999999
$ws.Cells.Item(12,2).Value = 46
$ws.Cells.Item(12,3).Value = -3.460514268753922
$ws.Cells.Item(12,4).Value = 97.26632312083132
$ws.Cells.Item(13,1).Value = 45032.99999999999
$ws.Cells.Item(13,2).Value = 46
$ws.Cells.Item(13,3).Value = 1.020549007077899
$ws.Cells.Item(13,4).Value = 94.8397173535453
$ws.Cells.Item(14,1).Value = 45039.99999999999
$ws.Cells.Item(14,2).Value = 47
$ws.Cells.Item(14,3).Value = -1.511815378292812
$ws.Cells.Item(14,4).Value = 91.70797270005116
$ws.Cells.Item(15,1).Value = 45046.99999999999
$ws.Cells.Item(15,2).Value = 47
$ws.Cells.Item(15,3).Value = 4.167649262036482
$ws.Cells.Item(15,4).Value = 94.45702813251414
$ws.Cells.Item(16,1).Value = 45060.99999999999
$ws.Cells.Item(16,2).Value = 48
$ws.Cells.Item(16,3).Value = 3.796646682899737
$ws.Cells.Item(16,4).Value = 93.86513260102345
$ws.Cells.Item(17,1).Value = 45067.99999999999
$ws.Cells.Item(17,2).Value = 49
$ws.Cells.Item(17,3).Value = 1.065597253771691
$ws.Cells.Item(17,4).Value = 94.63200973089894
$ws.Cells.Item(18,1).Value = 45074.99999999999
$ws.Cells.Item(18,2).Value = 49
$ws.Cells.Item(18,3).Value = 1.516866285569994
$ws.Cells.Item(18,4).Value = 96.53043661157996
$ws.Cells.Item(19,1).Value = 45088.99999999999
$ws.Cells.Item(19,2).Value = 50
$ws.Cells.Item(19,3).Value = 6.727914420272598
$ws.Cells.Item(19,4).Value = 96.94866099837077
$ws.Cells.Item(20,1).Value = 45095.99999999999
$ws.Cells.Item(20,2).Value = 51
$ws.Cells.Item(20,3).Value = 6.037476670497309
$ws.Cells.Item(20,4).Value = 96.34701708825753
$ws.Cells.Item(21,1).Value = 45102.99999999999
$ws.Cells.Item(21,2).Value = 51
$ws.Cells.Item(21,3).Value = 5.744315752500687
$ws.Cells.Item(21,4).Value = 100.0972925381883
$ws.Cells.Item(22,1).Value = 45109.99999999999
$ws.Cells.Item(22,2).Value = 52
$ws.Cells.Item(22,3).Value = 4.247183451549431
$ws.Cells.Item(22,4).Value = 97.75464153381289
$ws.Cells.Item(23,1).Value = 45116.99999999999
$ws.Cells.Item(23,2).Value = 52
$ws.Cells.Item(23,3).Value = 5.588196752390227
$ws.Cells.Item(23,4).Value = 99.46312860694228
$ws.Cells.Item(24,1).Value = 45123.99999999999
$ws.Cells.Item(24,2).Value = 53
$ws.Cells.Item(24,3).Value = 5.859388511042437
$ws.Cells.Item(24,4).Value = 97.31837504362758
$ws.Cells.Item(25,1).Value = 45130.99999999999
$ws.Cells.Item(25,2).Value = 53
$ws.Cells.Item(25,3).Value = 7.454478839641248
$ws.Cells.Item(25,4).Value = 98.21571782405998
$ws.Cells.Item(26,1).Value = 45137.99999999999
$ws.Cells.Item(26,2).Value = 54
$ws.Cells.Item(26,3).Value = 5.533402466730107
$ws.Cells.Item(26,4).Value = 100.7317607389798
$ws.Cells.Item(27,1).Value = 45144.99999999999
$ws.Cells.Item(27,2).Value = 54
$ws.Cells.Item(27,3).Value = 5.725128476521015
$ws.Cells.Item(27,4).Value = 100.4474920663119
$ws.Cells.Item(28,1).Value = 45151.99999999999
$ws.Cells.Item(28,2).Value = 55
$ws.Cells.Item(28,3).Value = 9.363505862273739
$ws.Cells.Item(28,4).Value = 101.5612983353066
$ws.Cells.Item(29,1).Value = 45158.99999999999
$ws.Cells.Item(29,2).Value = 55
$ws.Cells.Item(29,3).Value = 11.8217259455214
$ws.Cells.Item(29,4).Value = 103.9779155963288
$ws.Cells.Item(30,1).Value = 45165.99999999999
$ws.Cells.Item(30,2).Value = 56
$ws.Cells.Item(30,3).Value = 7.086322749268242
$ws.Cells.Item(30,4).Value = 99.11764648869305

$null = $ws.Range("A1").Select()

# --- Leave the workbook's active sheet/selection as it was before the
#     edit (the diff doesn't touch bookViews/activeTab). ---
$null = $wsWeekly.Activate()
$null = $wsWeekly.Range("A1").Select()
